$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B5 should become a real numeric value (3) instead of a text "3"
$ws.Range("B5").Value = 3

# Add new row 6 with the new annotation data
$ws.Range("A6").Value = "Ruilin"
$ws.Range("B6").Value = "'4"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "achieves speedups"
$ws.Range("D6").Value = "APC"
$ws.Range("E6").Value = "RES"
$ws.Range("F6").Value = "42be9703-0e9b-4ce8-962d-60bf1f233ce8"
$ws.Range("G6").Value = "SJCPLLpaW_annotated.xlsx"
$ws.Range("H6").Value = "The results show that DeePa achieves speedups compared to PyTorch and TensorFlow with all of the tested minibatch sizes."
